$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$ws.Range("A7").Value = 130993250
$ws.Range("M7").Value = "äldre spår"
$ws.Range("P7").Value = "Stötetorpet, Stötetorpet, Boh"
$ws.Range("Z7").Value = "14:53"
$ws.Range("AB7").Value = "14:53"
